$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(47, 8).Value = 7600
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 7600
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 7600
$ws.Cells.Item(47, 13).Value = $null
$ws.Cells.Item(47, 14).Value = -9544
$ws.Cells.Item(121, 8).Value = 1330.3334
$ws.Cells.Item(121, 10).Value = 1470.5
$ws.Cells.Item(121, 12).Value = 4411.5
$ws.Cells.Item(121, 14).Value = -7905.5
$ws.Cells.Item(125, 8).Value = 62501240
$ws.Cells.Item(125, 9).Value = 111111864
$ws.Cells.Item(125, 10).Value = 1867.1428
$ws.Cells.Item(125, 11).Value = 1000006776
$ws.Cells.Item(125, 12).Value = 16804.2852
$ws.Cells.Item(125, 13).Value = -1000004316
$ws.Cells.Item(125, 14).Value = -21724.2852
$ws.Cells.Item(132, 8).Value = 1698.4348
$ws.Cells.Item(132, 9).Value = 785.8378
$ws.Cells.Item(132, 10).Value = 5450.222
$ws.Cells.Item(132, 11).Value = 2357.5134
$ws.Cells.Item(132, 12).Value = 16350.666
$ws.Cells.Item(132, 13).Value = 172.4866000000002
$ws.Cells.Item(132, 14).Value = -21410.666
$ws.Cells.Item(137, 8).Value = 842.375
$ws.Cells.Item(137, 9).Value = 755.9286
$ws.Cells.Item(137, 10).Value = 909.6111
$ws.Cells.Item(137, 11).Value = 2267.7858
$ws.Cells.Item(137, 12).Value = 2728.8333
$ws.Cells.Item(137, 13).Value = 282.2142000000003
$ws.Cells.Item(137, 14).Value = -7828.8333
$ws.Cells.Item(138, 8).Value = 2629.3281
$ws.Cells.Item(138, 9).Value = 1205.9318
$ws.Cells.Item(138, 10).Value = 5760.8
$ws.Cells.Item(138, 11).Value = 3617.7954
$ws.Cells.Item(138, 12).Value = 17282.4
$ws.Cells.Item(138, 13).Value = 1522.2046
$ws.Cells.Item(138, 14).Value = -27562.4
$ws.Cells.Item(141, 8).Value = 2132.4856
$ws.Cells.Item(141, 9).Value = 2132.4856
$ws.Cells.Item(141, 11).Value = 6397.4568
$ws.Cells.Item(141, 13).Value = -1217.4568

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1070.1818
$ws.Cells.Item(61, 9).Value = 690.5454999999999
$ws.Cells.Item(61, 10).Value = 2209.0908
$ws.Cells.Item(61, 11).Value = 690.5454999999999
$ws.Cells.Item(61, 12).Value = 2209.0908
$ws.Cells.Item(61, 13).Value = -478.5454999999999
$ws.Cells.Item(61, 14).Value = -2633.0908
$ws.Cells.Item(74, 8).Value = 204910.52
$ws.Cells.Item(74, 9).Value = 228090.94
$ws.Cells.Item(74, 10).Value = 922.8
$ws.Cells.Item(74, 11).Value = 228090.94
$ws.Cells.Item(74, 12).Value = 922.8
$ws.Cells.Item(74, 13).Value = -227216.94
$ws.Cells.Item(74, 14).Value = -2670.8
$ws.Cells.Item(77, 8).Value = 204910.52
$ws.Cells.Item(77, 9).Value = 228090.94
$ws.Cells.Item(77, 10).Value = 922.8
$ws.Cells.Item(77, 11).Value = 1140454.7
$ws.Cells.Item(77, 12).Value = 4614
$ws.Cells.Item(77, 13).Value = -1136086.7
$ws.Cells.Item(77, 14).Value = -13350
$ws.Cells.Item(132, 8).Value = 1059.8431
$ws.Cells.Item(132, 9).Value = 833.3022999999999
$ws.Cells.Item(132, 10).Value = 2277.5
$ws.Cells.Item(132, 11).Value = 2499.9069
$ws.Cells.Item(132, 12).Value = 6832.5
$ws.Cells.Item(132, 13).Value = 30.09310000000005
$ws.Cells.Item(132, 14).Value = -11892.5
$ws.Cells.Item(136, 8).Value = 1070.1818
$ws.Cells.Item(136, 9).Value = 690.5454999999999
$ws.Cells.Item(136, 10).Value = 2209.0908
$ws.Cells.Item(136, 11).Value = 2071.6365
$ws.Cells.Item(136, 12).Value = 6627.2724
$ws.Cells.Item(136, 13).Value = 478.3635000000004
$ws.Cells.Item(136, 14).Value = -11727.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 806.1818
$ws.Cells.Item(134, 9).Value = 615.0217
$ws.Cells.Item(134, 10).Value = 1783.2222
$ws.Cells.Item(134, 11).Value = 1845.0651
$ws.Cells.Item(134, 12).Value = 5349.6666
$ws.Cells.Item(134, 13).Value = 689.9349
$ws.Cells.Item(134, 14).Value = -10419.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7577588.5
$ws.Cells.Item(31, 9).Value = 9092316
$ws.Cells.Item(31, 10).Value = 3954
$ws.Cells.Item(31, 11).Value = 9092316
$ws.Cells.Item(31, 12).Value = 3954
$ws.Cells.Item(31, 13).Value = -9092021
$ws.Cells.Item(31, 14).Value = -4544
$ws.Cells.Item(34, 8).Value = 7577588.5
$ws.Cells.Item(34, 9).Value = 9092316
$ws.Cells.Item(34, 10).Value = 3954
$ws.Cells.Item(34, 11).Value = 9092316
$ws.Cells.Item(34, 12).Value = 3954
$ws.Cells.Item(34, 13).Value = -9092114
$ws.Cells.Item(34, 14).Value = -4358
$ws.Cells.Item(58, 8).Value = 967.4888999999999
$ws.Cells.Item(58, 9).Value = 832.14813
$ws.Cells.Item(58, 10).Value = 1170.5
$ws.Cells.Item(58, 11).Value = 832.14813
$ws.Cells.Item(58, 12).Value = 1170.5
$ws.Cells.Item(58, 13).Value = -629.14813
$ws.Cells.Item(58, 14).Value = -1576.5
$ws.Cells.Item(112, 8).Value = 28000
$ws.Cells.Item(112, 10).Value = 28000
$ws.Cells.Item(112, 12).Value = 28000
$ws.Cells.Item(112, 14).Value = -30954
$ws.Cells.Item(132, 8).Value = 1761.8462
$ws.Cells.Item(132, 9).Value = 1545.7894
$ws.Cells.Item(132, 10).Value = 2348.2856
$ws.Cells.Item(132, 11).Value = 4637.3682
$ws.Cells.Item(132, 12).Value = 7044.8568
$ws.Cells.Item(132, 13).Value = -2107.3682
$ws.Cells.Item(132, 14).Value = -12104.8568
$ws.Cells.Item(134, 8).Value = 1081.3572
$ws.Cells.Item(134, 9).Value = 1097.919
$ws.Cells.Item(134, 10).Value = 958.8
$ws.Cells.Item(134, 11).Value = 3293.757000000001
$ws.Cells.Item(134, 12).Value = 2876.4
$ws.Cells.Item(134, 13).Value = -758.7570000000005
$ws.Cells.Item(134, 14).Value = -7946.4
$ws.Cells.Item(136, 8).Value = 967.4888999999999
$ws.Cells.Item(136, 9).Value = 832.14813
$ws.Cells.Item(136, 10).Value = 1170.5
$ws.Cells.Item(136, 11).Value = 2496.44439
$ws.Cells.Item(136, 12).Value = 3511.5
$ws.Cells.Item(136, 13).Value = 53.55560999999989
$ws.Cells.Item(136, 14).Value = -8611.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value = 1149
$ws.Cells.Item(47, 9).Value = 298
$ws.Cells.Item(47, 10).Value = 2000
$ws.Cells.Item(47, 11).Value = 894
$ws.Cells.Item(47, 12).Value = 6000
$ws.Cells.Item(47, 13).Value = -463
$ws.Cells.Item(47, 14).Value = -6862
$ws.Cells.Item(48, 8).Value = 2142
$ws.Cells.Item(48, 9).Value = 1000
$ws.Cells.Item(48, 10).Value = 3055.6
$ws.Cells.Item(48, 11).Value = 3000
$ws.Cells.Item(48, 12).Value = 9166.799999999999
$ws.Cells.Item(48, 13).Value = -2750
$ws.Cells.Item(48, 14).Value = -9666.799999999999
$ws.Cells.Item(49, 8).Value = 2985
$ws.Cells.Item(49, 10).Value = 2985
$ws.Cells.Item(49, 12).Value = 8955
$ws.Cells.Item(49, 14).Value = -9267
$ws.Cells.Item(131, 8).Value = 7353792.5
$ws.Cells.Item(131, 10).Value = 9434844
$ws.Cells.Item(131, 12).Value = 28304532
$ws.Cells.Item(131, 14).Value = -28314612
$ws.Cells.Item(137, 8).Value = 7480
$ws.Cells.Item(137, 9).Value = 17502.5
$ws.Cells.Item(137, 10).Value = 3471
$ws.Cells.Item(137, 11).Value = 52507.5
$ws.Cells.Item(137, 12).Value = 10413
$ws.Cells.Item(137, 13).Value = -47407.5
$ws.Cells.Item(137, 14).Value = -20613

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 18485.848
$ws.Cells.Item(132, 9).Value = 23284.762
$ws.Cells.Item(132, 10).Value = 1505.0769
$ws.Cells.Item(132, 11).Value = 69854.28599999999
$ws.Cells.Item(132, 12).Value = 4515.2307
$ws.Cells.Item(132, 13).Value = -67324.28599999999
$ws.Cells.Item(132, 14).Value = -9575.2307

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(47, 8).Value = 5399
$ws.Cells.Item(47, 10).Value = 5399
$ws.Cells.Item(47, 12).Value = 5399
$ws.Cells.Item(47, 14).Value = -6379
$ws.Cells.Item(48, 8).Value = 5800
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 13).Value = $null
$ws.Cells.Item(52, 8).Value = 5399
$ws.Cells.Item(52, 10).Value = 5399
$ws.Cells.Item(52, 12).Value = 5399
$ws.Cells.Item(52, 14).Value = -5865
$ws.Cells.Item(132, 8).Value = 3332.6
$ws.Cells.Item(132, 9).Value = 3753.9666
$ws.Cells.Item(132, 10).Value = 2068.5
$ws.Cells.Item(132, 11).Value = 11261.8998
$ws.Cells.Item(132, 12).Value = 6205.5
$ws.Cells.Item(132, 13).Value = -8731.899800000001
$ws.Cells.Item(132, 14).Value = -11265.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).Value = $null
$ws.Cells.Item(49, 8).Value = 4831.3335
$ws.Cells.Item(49, 9).Value = 3000
$ws.Cells.Item(49, 10).Value = 4997.8184
$ws.Cells.Item(49, 11).Value = 3000
$ws.Cells.Item(49, 12).Value = 4997.8184
$ws.Cells.Item(49, 13).Value = -2770
$ws.Cells.Item(49, 14).Value = -5457.8184
$ws.Cells.Item(132, 8).Value = 18940354
$ws.Cells.Item(132, 9).Value = 24039346
$ws.Cells.Item(132, 10).Value = 1240.3572
$ws.Cells.Item(132, 11).Value = 72118038
$ws.Cells.Item(132, 12).Value = 3721.0716
$ws.Cells.Item(132, 13).Value = -72115508
$ws.Cells.Item(132, 14).Value = -8781.071599999999
